$d = $word.ActiveDocument
$lastPara = $d.Paragraphs.Last
$r = $d.Range($lastPara.Range.Start, $d.Content.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="left"/><w:rPr><w:rFonts w:hint="default"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:highlight w:val="none"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:highlight w:val="none"/><w:lang w:val="ru-RU"/></w:rPr><w:t>А расстояние было перенесено для линков групп, потому что вот там оно как раз таки нужно.</w:t></w:r></w:p>    <w:p>
      <w:pPr>
        <w:pBdr>
          <w:bottom w:val="single" w:color="auto" w:sz="12" w:space="0"/>
        </w:pBdr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Function Link Welder :</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">Node_A  = </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
        <w:t>Корневой Нод</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">Node_B = </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
        <w:t>Нод ветвь</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Source = SV|AI</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
        <w:t>Сначала нужно проверить наличие в памяти словаря нашего корневого нода, потом надо проверить наличие ветви в этом словаре</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
        <w:t>(Заглавные и строчные это 2 разных нода, и их связи тоже разные)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
        <w:t>(Хотя можно было бы оптимизировать это и попробовать ссылаться из заглавного словаря на прописной, но со скаляром, хотя в этом логики мало так что пока что бракуем это)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
        <w:t>Источник будет влиять на итог в рассчетах, т.е несмотря на то что параметры у нас считаются отношением, конечный результат это не будет прямой кореляцией отношения, а еще и доверия к источнику этой информации, хотя конечно по началу едва ли это будет сильно влиять на что либо.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pBdr>
          <w:bottom w:val="single" w:color="auto" w:sz="12" w:space="0"/>
        </w:pBdr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
        <w:t>Проходя через нейрон или любой раз как мы будем ссылаться на него, мы должны увеличивать его заряд</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pBdr>
          <w:bottom w:val="single" w:color="auto" w:sz="12" w:space="0"/>
        </w:pBdr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
        <w:t>Формирование связи требует определенного уровня заряда, так что мы не можем выучить все слово за раз</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:highlight w:val="none"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
